$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The second table (old rows 16-25, cols A-C) is merged alongside the first
# table as columns E-G of rows 2-11; clear the old tail rows first.
$ws.Range("A16:G25").ClearContents()

# Rewrite the full A1:G15 block for the merged two-table layout.
# A leading "'" forces genuinely numeric-looking / empty text to stay text,
# matching the source data (everything here is inline text, not numbers).
$ws.Cells.Item(1,1).Value = "检验项目"
$ws.Cells.Item(1,2).Value = "结果"
$ws.Cells.Item(1,3).Value = "参考范围"
$ws.Cells.Item(1,4).Value = "单位"
$ws.Cells.Item(1,5).Value = "检验项目"
$ws.Cells.Item(1,6).Value = "结果"
$ws.Cells.Item(1,7).Value = "参考范围"
$ws.Cells.Item(2,1).Value = "白细胞"
$ws.Cells.Item(2,2).Value = "'2.76"
$ws.Cells.Item(2,3).Value = "4.00-10.0"
$ws.Cells.Item(2,4).Value = "10^9/L"
$ws.Cells.Item(2,5).Value = "嗜碱性粒细胞比率"
$ws.Cells.Item(2,6).Value = "'0.0"
$ws.Cells.Item(2,7).Value = "0.0-1.0"
$ws.Cells.Item(3,1).Value = "红细胞"
$ws.Cells.Item(3,2).Value = "'4.45"
$ws.Cells.Item(3,3).Value = "3.50-5.50"
$ws.Cells.Item(3,4).Value = "10^12/L"
$ws.Cells.Item(3,5).Value = "中性细胞数"
$ws.Cells.Item(3,6).Value = "'1.1"
$ws.Cells.Item(3,7).Value = "2.0-7.7109"
$ws.Cells.Item(4,1).Value = "血红蛋白"
$ws.Cells.Item(4,2).Value = "'122.0"
$ws.Cells.Item(4,3).Value = "010.0-160."
$ws.Cells.Item(4,4).Value = "g/L"
$ws.Cells.Item(4,5).Value = "淋巴细胞数"
$ws.Cells.Item(4,6).Value = "'1.39"
$ws.Cells.Item(4,7).Value = "0.80-4.00109"
$ws.Cells.Item(5,1).Value = "红细胞压积"
$ws.Cells.Item(5,2).Value = "'35.0"
$ws.Cells.Item(5,3).Value = "436.0-50.0"
$ws.Cells.Item(5,4).Value = "'"
$ws.Cells.Item(5,5).Value = "单核细胞数"
$ws.Cells.Item(5,6).Value = "'0.21"
$ws.Cells.Item(5,7).Value = "0.12-0.80109"
$ws.Cells.Item(6,1).Value = "血小板"
$ws.Cells.Item(6,2).Value = "'"
$ws.Cells.Item(6,3).Value = "'307"
$ws.Cells.Item(6,4).Value = "10^9/L"
$ws.Cells.Item(6,5).Value = "嗜酸性粒细胞数"
$ws.Cells.Item(6,6).Value = "'0.1"
$ws.Cells.Item(6,7).Value = "0.05-0.50109"
$ws.Cells.Item(7,1).Value = "'"
$ws.Cells.Item(7,2).Value = "'10.0"
$ws.Cells.Item(7,3).Value = "9.0-13.0"
$ws.Cells.Item(7,4).Value = "'"
$ws.Cells.Item(7,5).Value = "嗜碱性粒细胞"
$ws.Cells.Item(7,6).Value = "'0.0"
$ws.Cells.Item(7,7).Value = "0.00-0.10109"
$ws.Cells.Item(8,1).Value = "血小板压积"
$ws.Cells.Item(8,2).Value = "'0.3"
$ws.Cells.Item(8,3).Value = "'"
$ws.Cells.Item(8,4).Value = "21红细胞分布宽度"
$ws.Cells.Item(8,5).Value = "红细胞分布宽度"
$ws.Cells.Item(8,6).Value = "'38.8"
$ws.Cells.Item(8,7).Value = "37.0-50.0"
$ws.Cells.Item(9,1).Value = "红细胞平均体积"
$ws.Cells.Item(9,2).Value = "'79.0"
$ws.Cells.Item(9,3).Value = "6486.0-100."
$ws.Cells.Item(9,4).Value = "'"
$ws.Cells.Item(9,5).Value = "RDW-CV"
$ws.Cells.Item(9,6).Value = "'13.6"
$ws.Cells.Item(9,7).Value = "'"
$ws.Cells.Item(10,1).Value = "平均血红蛋白量"
$ws.Cells.Item(10,2).Value = "'27.4"
$ws.Cells.Item(10,3).Value = "26.0-33.0"
$ws.Cells.Item(10,4).Value = "23血小板分布宽度"
$ws.Cells.Item(10,5).Value = "血小板分布宽度"
$ws.Cells.Item(10,6).Value = "'10.8"
$ws.Cells.Item(10,7).Value = "9.0-17.0"
$ws.Cells.Item(11,1).Value = "平均血红蛋白浓度"
$ws.Cells.Item(11,2).Value = "'"
$ws.Cells.Item(11,3).Value = "310-370"
$ws.Cells.Item(11,4).Value = "24大型血小板比率"
$ws.Cells.Item(11,5).Value = "大型血小板比率"
$ws.Cells.Item(11,6).Value = "'24.0"
$ws.Cells.Item(11,7).Value = "13.0-43.0"
$ws.Cells.Item(12,1).Value = "中性细胞比率"
$ws.Cells.Item(12,2).Value = "'38.0"
$ws.Cells.Item(12,3).Value = "445.0-77.0"
$ws.Cells.Item(12,4).Value = "'"
$ws.Cells.Item(12,5).Value = "'"
$ws.Cells.Item(12,6).Value = "'"
$ws.Cells.Item(12,7).Value = "'"
$ws.Cells.Item(13,1).Value = "淋巴细胞比率"
$ws.Cells.Item(13,2).Value = "'50.4"
$ws.Cells.Item(13,3).Value = "20.0-40.0"
$ws.Cells.Item(13,4).Value = "'"
$ws.Cells.Item(13,5).Value = "'"
$ws.Cells.Item(13,6).Value = "'"
$ws.Cells.Item(13,7).Value = "'"
$ws.Cells.Item(14,1).Value = "单核细胞比率"
$ws.Cells.Item(14,2).Value = "'7.6"
$ws.Cells.Item(14,3).Value = "3.0-8.0"
$ws.Cells.Item(14,4).Value = "'"
$ws.Cells.Item(14,5).Value = "'"
$ws.Cells.Item(14,6).Value = "'"
$ws.Cells.Item(14,7).Value = "'"
$ws.Cells.Item(15,1).Value = "嗜酸性粒细胞比率"
$ws.Cells.Item(15,2).Value = "'"
$ws.Cells.Item(15,3).Value = "0.5-5.0"
$ws.Cells.Item(15,4).Value = "'"
$ws.Cells.Item(15,5).Value = "'"
$ws.Cells.Item(15,6).Value = "'"
$ws.Cells.Item(15,7).Value = "'"
